# Sara Alert import format: add "Vaccine 1" and "Vaccine 2" column groups
# (10 new columns, CY:DH) so the vaccine table can be populated on import.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column headers (row 1) ------------------------------------------------
$ws.Range("CY1").Value = "Vaccine 1 Group Name"
$ws.Range("CZ1").Value = "Vaccine 1 Product Name"
$ws.Range("DA1").NumberFormat = "@"
$ws.Range("DA1").Value = "Vaccine 1 Administration Date"
$ws.Range("DB1").Value = "Vaccine 1 Dose Number"
$ws.Range("DC1").Value = "Vaccine 1 Notes"
$ws.Range("DD1").Value = "Vaccine 2 Group Name"
$ws.Range("DE1").Value = "Vaccine 2 Product Name"
$ws.Range("DF1").NumberFormat = "@"
$ws.Range("DF1").Value = "Vaccine 2 Administration Date"
$ws.Range("DG1").Value = "Vaccine 2 Dose Number"
$ws.Range("DH1").Value = "Vaccine 2 Notes"

# --- Sample / invalid test data (row 2) ------------------------------------
$ws.Range("CY2").Value = "abc"
$ws.Range("CZ2").Value = 456
$ws.Range("DA2").NumberFormat = "@"
$ws.Range("DA2").Value = "2020-00-00"
$ws.Range("DB2").Value = 0
$ws.Range("DC2").Value = 'aW$#$#$!T@TFE'
$ws.Range("DD2").Value = '##@$#@!#'
$ws.Range("DE2").Value = '#@$T@%#$!%'
$ws.Range("DF2").NumberFormat = "@"
$ws.Range("DF2").Value = '#@$!@#%$'
$ws.Range("DG2").Value = "bvsadegr"
$ws.Range("DH2").Value = '#@$R'

# --- Sample / invalid test data (row 3) ------------------------------------
$ws.Range("CY3").Value = 123
$ws.Range("CZ3").Value = "def"
$ws.Range("DA3").NumberFormat = "@"
$ws.Range("DA3").Value = "00123"
$ws.Range("DB3").Value = "abc"
$ws.Range("DC3").Value = "32ntkgmavd"
$ws.Range("DD3").Value = [char]0x2026 + "....."
$ws.Range("DE3").Value = [char]0x2026 + "..."
$ws.Range("DF3").NumberFormat = "@"
$ws.Range("DF3").Value = "abadsdfwaeber"
$ws.Range("DG3").Value = '#@%!$@#'
$ws.Range("DH3").Value = "ab1234"
